$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To do List")

# Row 4 (Item 0.3, "Fix all hardcoded date variables"): mark Status Complete
$ws.Range("G4").Value = "Complete"

# Row 9 (Item 1.3, "What is logic for Qtrs and Months..."): mark Status Complete
$ws.Range("G9").Value = "Complete"

# Row 11 (Item 1.5, "Build in all BUs/segments for Core NA and Total NA view"):
#   Notes -> "Lawrence to rethink UI behavior", Status -> "Pending"
$ws.Range("E11").Value = "Lawrence to rethink UI behavior"
$ws.Range("G11").Value = "Pending"

# Row 13 (Item 1.7, "Incorporate both standard Cost account hierarchy..."):
#   Notes -> "Low priority"
$ws.Range("E13").Value = "Low priority"

# Row 16 (Item 2.1, "Add section to right of YTD..."): update Notes text
$ws.Range("E16").Value = "Replace the text objects in upper left with ""Financial Summary"".  Iteration 2"

# Match the author's saved cursor position/view state
$ws.Range("I19").Select() | Out-Null
